# Update the "want-to-go" attendance counters (column F) for a few events.
# These same events are listed on both the "展览" sheet and the "全部类型" sheet,
# so both copies need to be updated.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 2 -> event id=79051 : 334 -> 335
    $ws.Range("F2").Value = 335

    # Row 4 -> event id=80943 : 1389 -> 1399
    $ws.Range("F4").Value = 1399
}

# Row 5 on "展览" -> event id=81658 : 661 -> 668
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 668

# Row 6 on "全部类型" -> event id=81658 : 661 -> 668
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 668
